# Add a new "2022-Q4" sheet (with its fund-holding detail data) right after
# the "总计" summary sheet and before "2022-Q3", then update the "总计"
# summary sheet with a new top row for 2022-Q4 (shifting the old rows down
# by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert a new 2022-Q4 row at the
#    top of the data (row 2), pushing the existing quarter rows down by
#    one (2022-Q3 -> row3, 2022-Q2 -> row4, 2022-Q1 -> row5).
#    Write bottom-up so we never clobber a value before it has been
#    shifted down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.34
$summary.Range("A5").Font.Bold = $true
$summary.Range("A5").Borders.LineStyle = 1
$summary.Range("A5").HorizontalAlignment = -4108
$summary.Range("A5").VerticalAlignment = -4160

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 2.1

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 24
$summary.Range("D3").Value = 5.02

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 28
$summary.Range("D2").Value = 4.02

# ---------------------------------------------------------------------
# 2) Insert a brand-new worksheet named "2022-Q4" right before "2022-Q3"
#    (i.e. right after "总计") holding the per-fund holdings detail.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($q3)
$ws.Name = "2022-Q4"

# Header row (B1:H1) - bold, thin box border, centered/top aligned -
# matches the header styling used on every other sheet in the workbook.
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"
$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 3) Fund-holding detail rows for 2022-Q4 (rows 2..29). Column A (row
#    index), like on every other sheet, reuses the bold/bordered style.
#    Fund codes / percentages are written as text (quoted with a
#    leading apostrophe) to match the source data, which keeps these
#    numeric-looking values as strings (preserves leading zeros, exact
#    decimal text, etc.) instead of converting them to numbers.
# ---------------------------------------------------------------------
    # row 2
    $ws.Range("A2").Value = 0
    $ws.Range("B2").Value = "'290011"
    $ws.Range("C2").Value = "泰信中小盘精选混合"
    $ws.Range("D2").Value = "'14.75"
    $ws.Range("E2").Value = "'94.33"
    $ws.Range("F2").Value = "'6.27"
    $ws.Range("G2").Value = "'0.9248"
    $ws.Range("H2").Value = 10
    # row 3
    $ws.Range("A3").Value = 1
    $ws.Range("B3").Value = "'007490"
    $ws.Range("C3").Value = "南方信息创新混合A"
    $ws.Range("D3").Value = "'15.05"
    $ws.Range("E3").Value = "'91.51"
    $ws.Range("F3").Value = "'3.99"
    $ws.Range("G3").Value = "'0.6005"
    $ws.Range("H3").Value = 9
    # row 4
    $ws.Range("A4").Value = 2
    $ws.Range("B4").Value = "'012155"
    $ws.Range("C4").Value = "汇添富成长先锋六个月持有期混合A"
    $ws.Range("D4").Value = "'14.21"
    $ws.Range("E4").Value = "'89.10"
    $ws.Range("F4").Value = "'3.71"
    $ws.Range("G4").Value = "'0.5272"
    $ws.Range("H4").Value = 7
    # row 5
    $ws.Range("A5").Value = 3
    $ws.Range("B5").Value = "'010599"
    $ws.Range("C5").Value = "汇添富高质量成长30一年持有期混合A"
    $ws.Range("D5").Value = "'14.06"
    $ws.Range("E5").Value = "'84.46"
    $ws.Range("F5").Value = "'3.41"
    $ws.Range("G5").Value = "'0.4794"
    $ws.Range("H5").Value = 8
    # row 6
    $ws.Range("A6").Value = 4
    $ws.Range("B6").Value = "'012650"
    $ws.Range("C6").Value = "博时半导体主题混合A"
    $ws.Range("D6").Value = "'7.14"
    $ws.Range("E6").Value = "'93.53"
    $ws.Range("F6").Value = "'3.99"
    $ws.Range("G6").Value = "'0.2849"
    $ws.Range("H6").Value = 8
    # row 7
    $ws.Range("A7").Value = 5
    $ws.Range("B7").Value = "'012651"
    $ws.Range("C7").Value = "博时半导体主题混合C"
    $ws.Range("D7").Value = "'4.61"
    $ws.Range("E7").Value = "'93.53"
    $ws.Range("F7").Value = "'3.99"
    $ws.Range("G7").Value = "'0.1839"
    $ws.Range("H7").Value = 8
    # row 8
    $ws.Range("A8").Value = 6
    $ws.Range("B8").Value = "'002580"
    $ws.Range("C8").Value = "泰信鑫选灵活配置混合C"
    $ws.Range("D8").Value = "'1.62"
    $ws.Range("E8").Value = "'93.92"
    $ws.Range("F8").Value = "'9.27"
    $ws.Range("G8").Value = "'0.1502"
    $ws.Range("H8").Value = 7
    # row 9
    $ws.Range("A9").Value = 7
    $ws.Range("B9").Value = "'009715"
    $ws.Range("C9").Value = "汇添富策略增长灵活配置混合"
    $ws.Range("D9").Value = "'3.51"
    $ws.Range("E9").Value = "'88.61"
    $ws.Range("F9").Value = "'4.24"
    $ws.Range("G9").Value = "'0.1488"
    $ws.Range("H9").Value = 7
    # row 10
    $ws.Range("A10").Value = 8
    $ws.Range("B10").Value = "'217012"
    $ws.Range("C10").Value = "招商行业领先混合A"
    $ws.Range("D10").Value = "'2.52"
    $ws.Range("E10").Value = "'86.52"
    $ws.Range("F10").Value = "'4.51"
    $ws.Range("G10").Value = "'0.1137"
    $ws.Range("H10").Value = 8
    # row 11
    $ws.Range("A11").Value = 9
    $ws.Range("B11").Value = "'960019"
    $ws.Range("C11").Value = "招商行業領先混合型證券投資基金 H"
    $ws.Range("D11").Value = "'2.52"
    $ws.Range("E11").Value = "'86.52"
    $ws.Range("F11").Value = "'4.51"
    $ws.Range("G11").Value = "'0.1137"
    $ws.Range("H11").Value = 8
    # row 12
    $ws.Range("A12").Value = 10
    $ws.Range("B12").Value = "'007491"
    $ws.Range("C12").Value = "南方信息创新混合C"
    $ws.Range("D12").Value = "'2.40"
    $ws.Range("E12").Value = "'91.51"
    $ws.Range("F12").Value = "'3.99"
    $ws.Range("G12").Value = "'0.0958"
    $ws.Range("H12").Value = 9
    # row 13
    $ws.Range("A13").Value = 11
    $ws.Range("B13").Value = "'001970"
    $ws.Range("C13").Value = "泰信鑫选灵活配置混合A"
    $ws.Range("D13").Value = "'0.94"
    $ws.Range("E13").Value = "'93.92"
    $ws.Range("F13").Value = "'9.27"
    $ws.Range("G13").Value = "'0.0871"
    $ws.Range("H13").Value = 7
    # row 14
    $ws.Range("A14").Value = 12
    $ws.Range("B14").Value = "'009058"
    $ws.Range("C14").Value = "博时科技创新混合C"
    $ws.Range("D14").Value = "'4.19"
    $ws.Range("E14").Value = "'77.03"
    $ws.Range("F14").Value = "'2.03"
    $ws.Range("G14").Value = "'0.0851"
    $ws.Range("H14").Value = 9
    # row 15
    $ws.Range("A15").Value = 13
    $ws.Range("B15").Value = "'009057"
    $ws.Range("C15").Value = "博时科技创新混合A"
    $ws.Range("D15").Value = "'3.82"
    $ws.Range("E15").Value = "'77.03"
    $ws.Range("F15").Value = "'2.03"
    $ws.Range("G15").Value = "'0.0775"
    $ws.Range("H15").Value = 9
    # row 16
    $ws.Range("A16").Value = 14
    $ws.Range("B16").Value = "'014703"
    $ws.Range("C16").Value = "博时时代领航混合A"
    $ws.Range("D16").Value = "'1.14"
    $ws.Range("E16").Value = "'78.99"
    $ws.Range("F16").Value = "'2.39"
    $ws.Range("G16").Value = "'0.0272"
    $ws.Range("H16").Value = 10
    # row 17
    $ws.Range("A17").Value = 15
    $ws.Range("B17").Value = "'011259"
    $ws.Range("C17").Value = "汇添富高质量成长30一年持有期混合C"
    $ws.Range("D17").Value = "'0.64"
    $ws.Range("E17").Value = "'84.46"
    $ws.Range("F17").Value = "'3.41"
    $ws.Range("G17").Value = "'0.0218"
    $ws.Range("H17").Value = 8
    # row 18
    $ws.Range("A18").Value = 16
    $ws.Range("B18").Value = "'012779"
    $ws.Range("C18").Value = "博时移动互联主题混合A"
    $ws.Range("D18").Value = "'0.80"
    $ws.Range("E18").Value = "'84.06"
    $ws.Range("F18").Value = "'2.54"
    $ws.Range("G18").Value = "'0.0203"
    $ws.Range("H18").Value = 8
    # row 19
    $ws.Range("A19").Value = 17
    $ws.Range("B19").Value = "'012156"
    $ws.Range("C19").Value = "汇添富成长先锋六个月持有期混合C"
    $ws.Range("D19").Value = "'0.40"
    $ws.Range("E19").Value = "'89.10"
    $ws.Range("F19").Value = "'3.71"
    $ws.Range("G19").Value = "'0.0148"
    $ws.Range("H19").Value = 7
    # row 20
    $ws.Range("A20").Value = 18
    $ws.Range("B20").Value = "'016238"
    $ws.Range("C20").Value = "华夏数字经济龙头混合C"
    $ws.Range("D20").Value = "'0.33"
    $ws.Range("E20").Value = "'91.36"
    $ws.Range("F20").Value = "'3.92"
    $ws.Range("G20").Value = "'0.0129"
    $ws.Range("H20").Value = 9
    # row 21
    $ws.Range("A21").Value = 19
    $ws.Range("B21").Value = "'010307"
    $ws.Range("C21").Value = "西藏东财信息产业精选混合A"
    $ws.Range("D21").Value = "'0.18"
    $ws.Range("E21").Value = "'87.60"
    $ws.Range("F21").Value = "'4.39"
    $ws.Range("G21").Value = "'0.0079"
    $ws.Range("H21").Value = 3
    # row 22
    $ws.Range("A22").Value = 20
    $ws.Range("B22").Value = "'016237"
    $ws.Range("C22").Value = "华夏数字经济龙头混合A"
    $ws.Range("D22").Value = "'0.20"
    $ws.Range("E22").Value = "'91.36"
    $ws.Range("F22").Value = "'3.92"
    $ws.Range("G22").Value = "'0.0078"
    $ws.Range("H22").Value = 9
    # row 23
    $ws.Range("A23").Value = 21
    $ws.Range("B23").Value = "'007439"
    $ws.Range("C23").Value = "东海科技动力混合A"
    $ws.Range("D23").Value = "'0.19"
    $ws.Range("E23").Value = "'83.30"
    $ws.Range("F23").Value = "'3.67"
    $ws.Range("G23").Value = "'0.0070"
    $ws.Range("H23").Value = 7
    # row 24
    $ws.Range("A24").Value = 22
    $ws.Range("B24").Value = "'014704"
    $ws.Range("C24").Value = "博时时代领航混合C"
    $ws.Range("D24").Value = "'0.25"
    $ws.Range("E24").Value = "'78.99"
    $ws.Range("F24").Value = "'2.39"
    $ws.Range("G24").Value = "'0.0060"
    $ws.Range("H24").Value = 10
    # row 25
    $ws.Range("A25").Value = 23
    $ws.Range("B25").Value = "'007463"
    $ws.Range("C25").Value = "东海科技动力混合C"
    $ws.Range("D25").Value = "'0.16"
    $ws.Range("E25").Value = "'83.30"
    $ws.Range("F25").Value = "'3.67"
    $ws.Range("G25").Value = "'0.0059"
    $ws.Range("H25").Value = 7
    # row 26
    $ws.Range("A26").Value = 24
    $ws.Range("B26").Value = "'005629"
    $ws.Range("C26").Value = "汇安趋势动力股票C"
    $ws.Range("D26").Value = "'0.08"
    $ws.Range("E26").Value = "'91.16"
    $ws.Range("F26").Value = "'5.96"
    $ws.Range("G26").Value = "'0.0048"
    $ws.Range("H26").Value = 8
    # row 27
    $ws.Range("A27").Value = 25
    $ws.Range("B27").Value = "'010308"
    $ws.Range("C27").Value = "西藏东财信息产业精选混合C"
    $ws.Range("D27").Value = "'0.09"
    $ws.Range("E27").Value = "'87.60"
    $ws.Range("F27").Value = "'4.39"
    $ws.Range("G27").Value = "'0.0040"
    $ws.Range("H27").Value = 3
    # row 28
    $ws.Range("A28").Value = 26
    $ws.Range("B28").Value = "'012780"
    $ws.Range("C28").Value = "博时移动互联主题混合C"
    $ws.Range("D28").Value = "'0.14"
    $ws.Range("E28").Value = "'84.06"
    $ws.Range("F28").Value = "'2.54"
    $ws.Range("G28").Value = "'0.0036"
    $ws.Range("H28").Value = 8
    # row 29
    $ws.Range("A29").Value = 27
    $ws.Range("B29").Value = "'005628"
    $ws.Range("C29").Value = "汇安趋势动力股票A"
    $ws.Range("D29").Value = "'0.02"
    $ws.Range("E29").Value = "'91.16"
    $ws.Range("F29").Value = "'5.96"
    $ws.Range("G29").Value = "'0.0012"
    $ws.Range("H29").Value = 8

$aColRange = $ws.Range("A2:A29")
$aColRange.Font.Bold = $true
$aColRange.Borders.LineStyle = 1
$aColRange.HorizontalAlignment = -4108
$aColRange.VerticalAlignment = -4160

Write-Output "2022-Q4 sheet added; 总计 summary updated"
